$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to reflect new date
$ws.Name = "Through 2022-10-29"

# Update header label text
$ws.Range("B1").Value = "October 2022 (through October 29)"

# Update existing cell values (incremented counts)
$ws.Range("BJ2").Value = 6
$ws.Range("B3").Value = 6
$ws.Range("B5").Value = 3
$ws.Range("L6").Value = 14
$ws.Range("BJ6").Value = 4
$ws.Range("B7").Value = 8
$ws.Range("V7").Value = 8
$ws.Range("BJ7").Value = 5
$ws.Range("B10").Value = 8
$ws.Range("L10").Value = 4
$ws.Range("AZ13").Value = 3
$ws.Range("AP15").Value = 2
$ws.Range("AF21").Value = 2
$ws.Range("B24").Value = 9
$ws.Range("L26").Value = 3
$ws.Range("BJ28").Value = 2
$ws.Range("K30").Value = 5
$ws.Range("L36").Value = 2
$ws.Range("V44").Value = 3
$ws.Range("B51").Value = 2
$ws.Range("L56").Value = 4
$ws.Range("B79").Value = 3
$ws.Range("AZ95").Value = 3
$ws.Range("B98").Value = 3

# Add new cell values for newly-reported incidents
$ws.Range("BJ9").Value = 1
$ws.Range("BT14").Value = 1
$ws.Range("V17").Value = 1
$ws.Range("B32").Value = 1
$ws.Range("AZ35").Value = 1
$ws.Range("L62").Value = 1
$ws.Range("BJ70").Value = 1
$ws.Range("B96").Value = 1
$ws.Range("B97").Value = 1
